$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 52: new pub-pedal entry (Farnah Green / The Bluebell)
$ws.Range("A52").Value = 43663
$ws.Range("B52").Value = "The Bluebell"
$ws.Range("C52").Value = "Farnah Green"
$ws.Range("D52").Value = "start/end at pub"
$ws.Range("E52").Value = 3.94
$ws.Range("F52").Value = 0.045902777777777772
$ws.Range("G52").Formula = "=F52/E52"
$ws.Range("H52").Value = 1
$ws.Range("I52").Value = 1
$ws.Range("J52").Value = 1
$ws.Range("N52").Value = 1
$ws.Range("O52").Value = "Pub pedal."
$ws.Range("P52").Formula = "=SUM(H52:N52)*E52"

# Row 54 totals: extend CL sum down to the new row
$ws.Range("H54").Formula = "=SUM(H5:H52)"

# Move the active selection, as recorded for the sheet view
$ws.Range("A52").Select()
